$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "77ab58b0-f669-4e6e-b26e-a8357789d080"
$ws.Range("A3").Value = "a2ec7b01-97e7-4b31-af0a-22843e0d5bee"
$ws.Range("A4").Value = "7052c8df-2109-410c-abdc-73e3daa3ba50"
$ws.Range("A5").Value = "e2cbd333-6275-4a76-bbad-aa32bea23534"
$ws.Range("A6").Value = "543110a9-c7fb-4592-a759-ecd8bde4dc23"
$ws.Range("A7").Value = "42ef002b-b4ac-40f9-910b-b83fc5bb7ef5"
$ws.Range("A8").Value = "41477a7c-27ad-4394-b514-bc99c3388894"
$ws.Range("A9").Value = "54832180-8534-4b92-9a57-4176f97aacfa"
$ws.Range("A10").Value = "6c9080bb-8061-4711-b02a-7da01338beb8"
$ws.Range("A11").Value = "c45d0f6d-e5c1-4d33-94f5-2e886b60bf19"
$ws.Range("A12").Value = "10ffe120-07be-4ac0-a3ee-2417c74c09ad"
$ws.Range("A13").Value = "4c42a4de-c483-4b00-9621-9bbfd1c3a6e5"
$ws.Range("A14").Value = "4e2f94bd-6e7f-47c6-99ad-48b96ff1a18b"
$ws.Range("A15").Value = "f5e1e006-f68b-4a12-a7ea-10f71a9f73a6"
$ws.Range("A16").Value = "40c98c1d-9f79-456e-ae83-a5703fa69e74"
$ws.Range("A17").Value = "da62cb63-8bf7-4b28-8536-bfa48dcdebb7"
$ws.Range("A18").Value = "78469df7-4600-4f34-a464-3abf5b3946b9"
$ws.Range("A19").Value = "dd901bee-2756-43e2-b2b0-0dfca5c04baa"
